# Update NATMI LR-pair values (Fn1-Itgb3) with refreshed TPM-derived numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values keyed by row -> column -> value
$updates = @{
    2  = @{ G = 6.206015333333333;  H = 18.618046;  I = 0.0150172404156507;  J = 0.0150172404156507;
            M = 0.110552;            N = 0.331656;   O = 0.01126249561724847; P = 0.01126249561724847;
            Q = 0.6860874071306666;  R = 6.174786664176; S = 0.0001691316043644326; T = 0.0001691316043644326 }
    3  = @{ G = 6.206015333333333;  H = 18.618046;  I = 0.0150172404156507;  J = 0.0150172404156507;
            O = 0.9181055646724333; P = 0.9181055646724334;
            Q = 55.92904874241689;  R = 503.361438681752; S = 0.01378741199163267; T = 0.01378741199163267 }
    4  = @{ G = 6.206015333333333;  H = 18.618046;  I = 0.0150172404156507;  J = 0.0150172404156507;
            M = 0.6933189999999999; N = 2.079957;    O = 0.07063193971031816; P = 0.07063193971031817;
            Q = 4.302748344891333;  R = 38.724735104022; S = 0.001060696819653593; T = 0.001060696819653593 }
    5  = @{ I = 0.9317452840597572; J = 0.9317452840597571;
            M = 0.110552;            N = 0.331656;   O = 0.01126249561724847; P = 0.01126249561724847;
            Q = 42.56832070029067;  R = 383.1148863026161; S = 0.01049377717811495; T = 0.01049377717811495 }
    6  = @{ I = 0.9317452840597572; J = 0.9317452840597571;
            O = 0.9181055646724333; P = 0.9181055646724334;
            S = 0.8554405301525602; T = 0.8554405301525602 }
    7  = @{ I = 0.9317452840597572; J = 0.9317452840597571;
            M = 0.6933189999999999; N = 2.079957;    O = 0.07063193971031816; P = 0.07063193971031817;
            Q = 266.9641936790363;  R = 2402.677743111327; S = 0.06581097672908204; T = 0.06581097672908204 }
    8  = @{ G = 22.00088566666667;  H = 66.002657;  I = 0.05323747552459213; J = 0.05323747552459213;
            M = 0.110552;            N = 0.331656;   O = 0.01126249561724847; P = 0.01126249561724847;
            Q = 2.432241912221333;  R = 21.890177209992; S = 0.0005995868347690917; T = 0.0005995868347690917 }
    9  = @{ G = 22.00088566666667;  H = 66.002657;  I = 0.05323747552459213; J = 0.05323747552459213;
            O = 0.9181055646724333; P = 0.9181055646724334;
            Q = 198.2735363572538; R = 1784.461827215284; S = 0.04887762252824051; T = 0.04887762252824051 }
    10 = @{ G = 22.00088566666667;  H = 66.002657;  I = 0.05323747552459213; J = 0.05323747552459213;
            M = 0.6933189999999999; N = 2.079957;    O = 0.07063193971031816; P = 0.07063193971031817;
            Q = 15.25363204952766;  R = 137.282688445749; S = 0.00376026616158253; T = 0.00376026616158253 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
